$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "73-32=41"
$t.Cell(1,2).Range.Text = "78-28=50"
$t.Cell(1,3).Range.Text = "30+51=81"
$t.Cell(1,4).Range.Text = "58+19=77"
$t.Cell(1,5).Range.Text = "65-12=53"
$t.Cell(2,1).Range.Text = "76-26=50"
$t.Cell(2,2).Range.Text = "27+13=40"
$t.Cell(2,3).Range.Text = "69-49=20"
$t.Cell(2,4).Range.Text = "77-43=34"
$t.Cell(2,5).Range.Text = "35+23=58"
$t.Cell(3,1).Range.Text = "12+45=57"
$t.Cell(3,2).Range.Text = "45+26=71"
$t.Cell(3,3).Range.Text = "24-2=22"
$t.Cell(3,4).Range.Text = "23+6=29"
$t.Cell(3,5).Range.Text = "63-61=2"
$t.Cell(4,1).Range.Text = "35-12=23"
$t.Cell(4,2).Range.Text = "74-41=33"
$t.Cell(4,3).Range.Text = "51-49=2"
$t.Cell(4,4).Range.Text = "91-77=14"
$t.Cell(4,5).Range.Text = "22+21=43"
$t.Cell(5,1).Range.Text = "12+63=75"
$t.Cell(5,2).Range.Text = "32+8=40"
$t.Cell(5,3).Range.Text = "11+73=84"
$t.Cell(5,4).Range.Text = "25+21=46"
$t.Cell(5,5).Range.Text = "41+16=57"
$t.Cell(6,1).Range.Text = "65-44=21"
$t.Cell(6,2).Range.Text = "75-25=50"
$t.Cell(6,3).Range.Text = "51+33=84"
$t.Cell(6,4).Range.Text = "75-69=6"
$t.Cell(6,5).Range.Text = "63+34=97"
$t.Cell(7,1).Range.Text = "32-16=16"
$t.Cell(7,2).Range.Text = "7+56=63"
$t.Cell(7,3).Range.Text = "2+43=45"
$t.Cell(7,4).Range.Text = "75-64=11"
$t.Cell(7,5).Range.Text = "97-8=89"
$t.Cell(8,1).Range.Text = "48+12=60"
$t.Cell(8,2).Range.Text = "90-33=57"
$t.Cell(8,3).Range.Text = "93-66=27"
$t.Cell(8,4).Range.Text = "42-5=37"
$t.Cell(8,5).Range.Text = "81+11=92"
$t.Cell(9,1).Range.Text = "41-33=8"
$t.Cell(9,2).Range.Text = "55+15=70"
$t.Cell(9,3).Range.Text = "52+28=80"
$t.Cell(9,4).Range.Text = "61-33=28"
$t.Cell(9,5).Range.Text = "65+19=84"
$t.Cell(10,1).Range.Text = "77-52=25"
$t.Cell(10,2).Range.Text = "3+83=86"
$t.Cell(10,3).Range.Text = "9+27=36"
$t.Cell(10,4).Range.Text = "89-23=66"
$t.Cell(10,5).Range.Text = "72-16=56"
$t.Cell(11,1).Range.Text = "68-1=67"
$t.Cell(11,2).Range.Text = "97-71=26"
$t.Cell(11,3).Range.Text = "39+27=66"
$t.Cell(11,4).Range.Text = "1+29=30"
$t.Cell(11,5).Range.Text = "38+56=94"
$t.Cell(12,1).Range.Text = "39+59=98"
$t.Cell(12,2).Range.Text = "44-33=11"
$t.Cell(12,3).Range.Text = "46+47=93"
$t.Cell(12,4).Range.Text = "98-63=35"
$t.Cell(12,5).Range.Text = "56-6=50"
$t.Cell(13,1).Range.Text = "99-77=22"
$t.Cell(13,2).Range.Text = "3+71=74"
$t.Cell(13,3).Range.Text = "65+2=67"
$t.Cell(13,4).Range.Text = "88-13=75"
$t.Cell(13,5).Range.Text = "2+44=46"
$t.Cell(14,1).Range.Text = "56-40=16"
$t.Cell(14,2).Range.Text = "4+16=20"
$t.Cell(14,3).Range.Text = "27+22=49"
$t.Cell(14,4).Range.Text = "46+12=58"
$t.Cell(14,5).Range.Text = "78+8=86"
$t.Cell(15,1).Range.Text = "90-76=14"
$t.Cell(15,2).Range.Text = "77+11=88"
$t.Cell(15,3).Range.Text = "41-23=18"
$t.Cell(15,4).Range.Text = "85-61=24"
$t.Cell(15,5).Range.Text = "44-38=6"
$t.Cell(16,1).Range.Text = "2+24=26"
$t.Cell(16,2).Range.Text = "2+14=16"
$t.Cell(16,3).Range.Text = "83-1=82"
$t.Cell(16,4).Range.Text = "98-84=14"
$t.Cell(16,5).Range.Text = "84+15=99"
$t.Cell(17,1).Range.Text = "73+24=97"
$t.Cell(17,2).Range.Text = "2+11=13"
$t.Cell(17,3).Range.Text = "34-17=17"
$t.Cell(17,4).Range.Text = "52-5=47"
$t.Cell(17,5).Range.Text = "51-30=21"
$t.Cell(18,1).Range.Text = "76+14=90"
$t.Cell(18,2).Range.Text = "23+5=28"
$t.Cell(18,3).Range.Text = "38+1=39"
$t.Cell(18,4).Range.Text = "48+11=59"
$t.Cell(18,5).Range.Text = "6+6=12"
$t.Cell(19,1).Range.Text = "20+22=42"
$t.Cell(19,2).Range.Text = "92-27=65"
$t.Cell(19,3).Range.Text = "88-16=72"
$t.Cell(19,4).Range.Text = "85-12=73"
$t.Cell(19,5).Range.Text = "70-35=35"
$t.Cell(20,1).Range.Text = "81-61=20"
$t.Cell(20,2).Range.Text = "81-74=7"
$t.Cell(20,3).Range.Text = "2+40=42"
$t.Cell(20,4).Range.Text = "16+53=69"
$t.Cell(20,5).Range.Text = "38-18=20"

Write-Output "Replaced 100 cell values"
